# Chore: changed kinetics column names to use plural (kineticsDB compatibility)
#
# Renames several header cells on the "kinetics1" sheet so the ref-type /
# ref columns and the effector columns use plural naming, and updates the
# sheet's selection / scroll position to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kinetics1")

# --- Header renames (row 1) -------------------------------------------------

# "negative effector" / "positive effector" -> plural
$ws.Range("H1").Value = "negative effectors"
$ws.Range("I1").Value = "positive effectors"

# (J1 = "allosteric", K1 = "subunits" stay the same)

# "*_ref_type" / "*_ref" -> "*_refs_type" / "*_refs"
$ws.Range("L1").Value = "mechanism_refs_type"
$ws.Range("M1").Value = "mechanism_refs"
$ws.Range("N1").Value = "inhibitors_refs_type"
$ws.Range("O1").Value = "inhibitors_refs"
$ws.Range("P1").Value = "activators_refs_type"
$ws.Range("Q1").Value = "activators_refs"
$ws.Range("R1").Value = "negative_effectors_refs_type"
$ws.Range("S1").Value = "negative_effectors_refs"
$ws.Range("T1").Value = "positive_effectors_refs_type"
$ws.Range("U1").Value = "positive_effectors_refs"
$ws.Range("V1").Value = "subunits_refs_type"
$ws.Range("W1").Value = "subunits_refs"

# (X1 = "comments" stays the same)

# --- View state --------------------------------------------------------------

# Keep "kinetics1" the active/selected sheet (activeTab stays pointed here).
$ws.Activate()

# Scroll so column C is the left-most visible column, then move the
# selection to X1 (matches the saved sheetView/selection state).
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("X1").Select()

# Slightly widen the sheet-tabs / scrollbar split (tabRatio 990 -> 993).
$excel.ActiveWindow.TabRatio = 993 / 1650
